# Insert two new price-report rows ("Conconina(o)" and "Escarola") right
# after the existing row 660, pushing all subsequent rows down by two.
# (Weekly Fruit/Vegetable price update - new week's records added.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(661).Insert()
$ws.Rows.Item(661).Insert()

# New row 661: Conconina(o)
$ws.Cells.Item(661,1).Value  = 10
$ws.Cells.Item(661,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(661,3).Value  = "La Araucanía"
$ws.Cells.Item(661,4).Value  = 44491
$ws.Cells.Item(661,5).Value  = 9
$ws.Cells.Item(661,6).Value  = 100112033
$ws.Cells.Item(661,7).Value  = "Lechuga"
$ws.Cells.Item(661,8).Value  = "Conconina(o)"
$ws.Cells.Item(661,9).Value  = "Primera"
$ws.Cells.Item(661,10).Value = 365
$ws.Cells.Item(661,11).Value = 7000
$ws.Cells.Item(661,12).Value = 8000
$ws.Cells.Item(661,13).Value = 7575
$ws.Cells.Item(661,14).Value = "$/caja 10 unidades"
$ws.Cells.Item(661,15).Value = "Región Metropolitana"
$ws.Cells.Item(661,16).Value = 758
$ws.Cells.Item(661,17).Value = 10
$ws.Cells.Item(661,18).Value = "Hortaliza"

# New row 662: Escarola
$ws.Cells.Item(662,1).Value  = 10
$ws.Cells.Item(662,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(662,3).Value  = "La Araucanía"
$ws.Cells.Item(662,4).Value  = 44172
$ws.Cells.Item(662,5).Value  = 9
$ws.Cells.Item(662,6).Value  = 100112033
$ws.Cells.Item(662,7).Value  = "Lechuga"
$ws.Cells.Item(662,8).Value  = "Escarola"
$ws.Cells.Item(662,9).Value  = "Primera"
$ws.Cells.Item(662,10).Value = 650
$ws.Cells.Item(662,11).Value = 7000
$ws.Cells.Item(662,12).Value = 8000
$ws.Cells.Item(662,13).Value = 7462
$ws.Cells.Item(662,14).Value = "$/caja 15 unidades"
$ws.Cells.Item(662,15).Value = "Región Metropolitana"
$ws.Cells.Item(662,16).Value = 497
$ws.Cells.Item(662,17).Value = 15
$ws.Cells.Item(662,18).Value = "Hortaliza"
